# Scheduled-runner refresh of cached market/profit figures across the
# Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). These are
# plain cached values (no formulas in the workbook), so the update is
# applied as direct cell writes per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 30.307692
$ws.Range("I11").Value = 30.307692
$ws.Range("K11").Value = 30.307692
$ws.Range("M11").Value = 109.692308
$ws.Range("H42").Value = 401.66666
$ws.Range("I42").Value = 362.66666
$ws.Range("J42").Value = 421.16666
$ws.Range("K42").Value = 1087.99998
$ws.Range("L42").Value = 1263.49998
$ws.Range("M42").Value = -857.9999800000001
$ws.Range("N42").Value = -1723.49998
$ws.Range("H62").Value = 4631.3335
$ws.Range("I62").Value = 4631.3335
$ws.Range("K62").Value = 4631.3335
$ws.Range("M62").Value = -4007.3335
$ws.Range("H65").Value = 4631.3335
$ws.Range("I65").Value = 4631.3335
$ws.Range("K65").Value = 23156.6675
$ws.Range("M65").Value = -20036.6675
$ws.Range("H112").Value = 4411.923
$ws.Range("J112").Value = 4411.923
$ws.Range("L112").Value = 13235.769
$ws.Range("N112").Value = -15451.769
$ws.Range("H141").Value = 5480.4614
$ws.Range("I141").Value = 1708
$ws.Range("J141").Value = 50750
$ws.Range("K141").Value = 5124
$ws.Range("L141").Value = 152250
$ws.Range("M141").Value = 56
$ws.Range("N141").Value = -162610
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16481.564
$ws.Range("I32").Value = 20944
$ws.Range("J32").Value = 5153.846
$ws.Range("K32").Value = 20944
$ws.Range("L32").Value = 5153.846
$ws.Range("M32").Value = -20657
$ws.Range("N32").Value = -5727.846
$ws.Range("H74").Value = 1326.6666
$ws.Range("I74").Value = 1154.1666
$ws.Range("J74").Value = 2016.6666
$ws.Range("K74").Value = 1154.1666
$ws.Range("L74").Value = 2016.6666
$ws.Range("M74").Value = -280.1666
$ws.Range("N74").Value = -3764.6666
$ws.Range("H77").Value = 1326.6666
$ws.Range("I77").Value = 1154.1666
$ws.Range("J77").Value = 2016.6666
$ws.Range("K77").Value = 5770.833000000001
$ws.Range("L77").Value = 10083.333
$ws.Range("M77").Value = -1402.833000000001
$ws.Range("N77").Value = -18819.333
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H102").Value = 2818.3333
$ws.Range("I102").Value = 2982
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2982
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -1360
$ws.Range("N102").Value = -5244
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2470032.2
$ws.Range("I80").Value = 7409607
$ws.Range("J80").Value = 244.7
$ws.Range("K80").Value = 7409607
$ws.Range("L80").Value = 244.7
$ws.Range("M80").Value = -7408609
$ws.Range("N80").Value = -2240.7
$ws.Range("H83").Value = 2470032.2
$ws.Range("I83").Value = 7409607
$ws.Range("J83").Value = 244.7
$ws.Range("K83").Value = 37048035
$ws.Range("L83").Value = 1223.5
$ws.Range("M83").Value = -37043043
$ws.Range("N83").Value = -11207.5
$ws.Range("H134").Value = 2181.5217
$ws.Range("I134").Value = 2162.4707
$ws.Range("J134").Value = 2235.5
$ws.Range("K134").Value = 6487.4121
$ws.Range("L134").Value = 6706.5
$ws.Range("M134").Value = -3952.4121
$ws.Range("N134").Value = -11776.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2874.45
$ws.Range("I132").Value = 2499.3225
$ws.Range("J132").Value = 4166.5557
$ws.Range("K132").Value = 7497.967500000001
$ws.Range("L132").Value = 12499.6671
$ws.Range("M132").Value = -4967.967500000001
$ws.Range("N132").Value = -17559.6671
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 159482.12
$ws.Range("I68").Value = 244449.12
$ws.Range("J68").Value = 1134.5454
$ws.Range("K68").Value = 733347.36
$ws.Range("L68").Value = 3403.6362
$ws.Range("M68").Value = -732536.36
$ws.Range("N68").Value = -5025.6362
$ws.Range("H71").Value = 159482.12
$ws.Range("I71").Value = 244449.12
$ws.Range("J71").Value = 1134.5454
$ws.Range("K71").Value = 2200042.08
$ws.Range("L71").Value = 10210.9086
$ws.Range("M71").Value = -2195986.08
$ws.Range("N71").Value = -18322.9086
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3236.818
$ws.Range("I80").Value = 3057.8572
$ws.Range("J80").Value = 3550
$ws.Range("K80").Value = 3057.8572
$ws.Range("L80").Value = 3550
$ws.Range("M80").Value = -2059.8572
$ws.Range("N80").Value = -5546
$ws.Range("H83").Value = 3236.818
$ws.Range("I83").Value = 3057.8572
$ws.Range("J83").Value = 3550
$ws.Range("K83").Value = 15289.286
$ws.Range("L83").Value = 17750
$ws.Range("M83").Value = -10297.286
$ws.Range("N83").Value = -27734
$ws.Range("H93").Value = 35125
$ws.Range("J93").Value = 35125
$ws.Range("L93").Value = 35125
$ws.Range("N93").Value = -38869
$ws.Range("H102").Value = 3728.5
$ws.Range("I102").Value = 3577.6667
$ws.Range("K102").Value = 3577.6667
$ws.Range("M102").Value = -1955.6667
$ws.Range("H126").Value = 3215.8333
$ws.Range("I126").Value = 2821.111
$ws.Range("J126").Value = 4400
$ws.Range("K126").Value = 8463.332999999999
$ws.Range("L126").Value = 13200
$ws.Range("M126").Value = -5993.332999999999
$ws.Range("N126").Value = -18140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1500
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1139
$ws.Range("N82").Value = ""
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1500
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -252
$ws.Range("N85").Value = ""
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 23166.334
$ws.Range("I75").Value = 9509
$ws.Range("K75").Value = 9509
$ws.Range("M75").Value = -8573
$ws.Range("H78").Value = 23166.334
$ws.Range("I78").Value = 9509
$ws.Range("K78").Value = 28527
$ws.Range("M78").Value = -23847
$ws.Range("H136").Value = 1897.25
$ws.Range("I136").Value = 1851.4
$ws.Range("J136").Value = 2011.875
$ws.Range("K136").Value = 5554.200000000001
$ws.Range("L136").Value = 6035.625
$ws.Range("M136").Value = -3004.200000000001
$ws.Range("N136").Value = -11135.625
